# Generate Report for Handback
#
# The "6b07ccfe-699f-489d-aef5-2635bc7b185d.md" file has now been handed
# back for both locales, so its Status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is reported (the
# Overview roll-up sheet plus each per-locale detail sheet), and each
# locale's "Latest Handback DateTime" is stamped with the new handback
# timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: roll-up Status columns for the 6b07ccfe row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn detail sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("G3").Value = "2016-02-22 13:54:13"

# --- de-de detail sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("G3").Value = "2016-02-22 13:54:38"
